$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "negative"
$ws.Range("C5").Value = -5
$ws.Range("D5").Value = -30
$ws.Range("E5").Value = -30
$ws.Range("F5").Value = "PASS"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "integer"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = "PASS"

$ws.Range("F6").Select()
